$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Utdanning"
$ws.Range("A2").Value = "Barnehagelærer"
$ws.Range("E3").Value = 1372
$ws.Range("C9").Value = 45741
$ws.Range("E14").Value = 11947
$ws.Range("C19").Value = 48121
$ws.Range("E22").Value = 13305
$ws.Range("A23").Value = "Grunnskolelærer"
$ws.Range("C24").Value = 52369
$ws.Range("E26").Value = 1417
$ws.Range("C28").Value = 51766
$ws.Range("C30").Value = 50917
$ws.Range("C34").Value = 49716
$ws.Range("C36").Value = 49424
$ws.Range("E37").Value = 9880
$ws.Range("E39").Value = 9967
$ws.Range("C43").Value = 50777
$ws.Range("A44").Value = "Faglærer"
$ws.Range("C44").Value = 16556
$ws.Range("E44").Value = 1
$ws.Range("C45").Value = 16581
$ws.Range("E45").Value = 997
$ws.Range("C46").Value = 16636
$ws.Range("E46").Value = 2010
$ws.Range("C47").Value = 16655
$ws.Range("E47").Value = 3089
$ws.Range("C48").Value = 16666
$ws.Range("E48").Value = 4178
$ws.Range("E49").Value = 5282
$ws.Range("E50").Value = 6377
$ws.Range("C51").Value = 16637
$ws.Range("E51").Value = 7461
$ws.Range("C52").Value = 16616
$ws.Range("E52").Value = 8524
$ws.Range("C53").Value = 16590
$ws.Range("E53").Value = 9577
$ws.Range("E54").Value = 10587
$ws.Range("C55").Value = 16565
$ws.Range("E55").Value = 11554
$ws.Range("C56").Value = 16565
$ws.Range("C57").Value = 16575
$ws.Range("C58").Value = 16591
$ws.Range("E58").Value = 14286
$ws.Range("C60").Value = 16643
$ws.Range("C61").Value = 16689
$ws.Range("E61").Value = 16831
$ws.Range("C62").Value = 16744
$ws.Range("E63").Value = 18459
$ws.Range("C64").Value = 16866
$ws.Range("E64").Value = 19251
$ws.Range("E65").Value = 1
$ws.Range("E85").Value = 13
$ws.Range("C86").Value = 13249
$ws.Range("E86").Value = -1
$ws.Range("C87").Value = 13275
$ws.Range("E87").Value = -50
$ws.Range("C88").Value = 13336
$ws.Range("E88").Value = -146
$ws.Range("C89").Value = 13396
$ws.Range("E89").Value = -256
$ws.Range("C90").Value = 13442
$ws.Range("E90").Value = -351
$ws.Range("C91").Value = 13513
$ws.Range("E91").Value = -489
$ws.Range("C92").Value = 13578
$ws.Range("E92").Value = -634
$ws.Range("C93").Value = 13612
$ws.Range("E93").Value = -746
$ws.Range("C94").Value = 13570
$ws.Range("E94").Value = -807
$ws.Range("C95").Value = 13502
$ws.Range("E95").Value = -838
$ws.Range("C96").Value = 13420
$ws.Range("E96").Value = -864
$ws.Range("C97").Value = 13362
$ws.Range("E97").Value = -927
$ws.Range("C98").Value = 13318
$ws.Range("E98").Value = -997
$ws.Range("C99").Value = 13301
$ws.Range("E99").Value = -1110
$ws.Range("C100").Value = 13235
$ws.Range("E100").Value = -1156
$ws.Range("C101").Value = 13143
$ws.Range("E101").Value = -1182
$ws.Range("C102").Value = 13061
$ws.Range("E102").Value = -1209
$ws.Range("C103").Value = 13033
$ws.Range("E103").Value = -1285
$ws.Range("C104").Value = 13033
$ws.Range("E104").Value = -1399
$ws.Range("C105").Value = 13037
$ws.Range("E105").Value = -1496
$ws.Range("C106").Value = 13053
$ws.Range("E106").Value = -1611